$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New timesheet rows -----------------------------------------------
# NOTE: the two new description strings are written in this specific
# order (row 13 before row 12) so that the workbook's shared-string
# table receives them in the same order as the authoritative edit
# (index 13 = "Inventory update..." / index 14 = "Consumable item...").
$ws.Cells.Item(13, 3).Value = "Inventory update when (un)equiping + Spear + Basic combat text"
$ws.Cells.Item(12, 3).Value = "Consumable item + Health potion + Character panel show equiped stat boost"

$ws.Cells.Item(12, 1).Value = 43071
$ws.Cells.Item(12, 1).NumberFormat = "d-mmm"
$ws.Cells.Item(12, 2).Value = 0.072916666666666671
$ws.Cells.Item(12, 2).NumberFormat = "h:mm"

$ws.Cells.Item(13, 1).Value = 43071
$ws.Cells.Item(13, 1).NumberFormat = "d-mmm"
$ws.Cells.Item(13, 2).Value = 0.038194444444444441
$ws.Cells.Item(13, 2).NumberFormat = "h:mm"

# --- Column C got a lot wider to fit the longer descriptions ----------
$ws.Columns.Item(3).ColumnWidth = 76.25

# --- Selection follows the last entered cell ---------------------------
$ws.Range("C13").Select() | Out-Null
